$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows above the current header row (row 1).
# This pushes the existing text header (row 1) down to row 3, and all
# the data rows (old 2..51) down to new rows 4..53.
$ws.Rows.Item(1).Resize(2).Insert()

# The bold/bordered header style that used to live on row 1 now sits on
# row 3 (it moved down along with the cells). Copy that formatting up
# onto the new row 1, then strip it back off row 3 so row 3 becomes a
# plain, unstyled row again.
$ws.Range("A3:N3").Copy()
$ws.Range("A1:N1").PasteSpecial(-4122)
$ws.Range("A3:N3").ClearFormats()

# New row 1: a 0-based numeric column-index header, replacing the old
# text header that now lives on row 3.
$ws.Range("A1").Value2 = 0
$ws.Range("B1").Value2 = 1
$ws.Range("C1").Value2 = 2
$ws.Range("D1").Value2 = 3
$ws.Range("E1").Value2 = 4
$ws.Range("F1").Value2 = 5
$ws.Range("G1").Value2 = 6
$ws.Range("H1").Value2 = 7
$ws.Range("I1").Value2 = 8
$ws.Range("J1").Value2 = 9
$ws.Range("K1").Value2 = 10
$ws.Range("L1").Value2 = 11
$ws.Range("M1").Value2 = 12
$ws.Range("N1").Value2 = 13

# New row 2: a "Washer" label in column E, with every other cell in the
# row touched (no-op formatting) so the whole row A2:N2 is materialized
# as present-but-blank cells, matching the rest of the row layout.
$ws.Range("A2:N2").Font.Bold = $false
$ws.Range("E2").Value2 = "Washer"
